# Regenerate merged AHB files
# 1. Rename header row labels from *_old/_new suffixes to *_FV2210/_FV2304.
# 2. Turn the header + data range into a native Excel Table ("Table1").
# 3. Freeze the header row (split at row 1, top-left cell A2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @(
    "Segmentname_FV2210",
    "Segmentgruppe_FV2210",
    "Segment_FV2210",
    "Datenelement_FV2210",
    "Segment ID_FV2210",
    "Code_FV2210",
    "Qualifier_FV2210",
    "Beschreibung_FV2210",
    "Bedingungsausdruck_FV2210",
    "Bedingung_FV2210",
    "diff",
    "Segmentname_FV2304",
    "Segmentgruppe_FV2304",
    "Segment_FV2304",
    "Datenelement_FV2304",
    "Segment ID_FV2304",
    "Code_FV2304",
    "Qualifier_FV2304",
    "Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304",
    "Bedingung_FV2304"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Convert the used range into a native table (creates xl/tables/table1.xml
# and wires up the worksheet's <tableParts>).
$dataRange = $ws.Range("A1:U64")
$tbl = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# Freeze panes above row 2 (header row stays visible while scrolling).
$ws.Range("A2").Select() | Out-Null
[void]($excel.ActiveWindow.FreezePanes = $true)
